# Add a new "SQL case expression" entry to the sql_lib workbook.
# Mirrors the existing rows: column A = "SQL", column B = short title,
# column C = the actual SQL snippet (wrapped text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the exact SQL snippet text (preserving the original trailing
# spaces / line breaks from the source data).
$sqlText = "SELECT CUBK, CUNBR, CUSTAT,    " + "`n" `
  + "  case custat                  " + "`n" `
  + "  when 'R' then 'assigned'     " + "`n" `
  + "  when ' ' then 'n/a'          " + "`n" `
  + "  end                          " + "`n" `
  + "FROM zusrlib/cup00301          "

# New row goes right after the current last row (35).
$newRow = 36

# Column A repeats the "SQL" language label used throughout the sheet.
$ws.Cells.Item($newRow, 1).Value = "SQL"

# Write column C (the SQL text) before column B ("Case expression") so
# that the shared-string table picks up the two new strings in the same
# order as the source workbook (SQL text first, title second).
$ws.Cells.Item($newRow, 3).Value = $sqlText
$ws.Cells.Item($newRow, 3).WrapText = $true

$ws.Cells.Item($newRow, 2).Value = "Case expression"
$ws.Cells.Item($newRow, 2).WrapText = $true

# Keep the row height consistent with the rest of the table.
$ws.Rows.Item($newRow).RowHeight = 37.5

# Update the selection to match the post-edit state (cursor moved past
# the newly-added row).
$ws.Range("A37").Select() | Out-Null
